# Corrected some selection scopes
# A new quarter row ("2020-04-01") was missing from the matched-errors
# series; insert it as row 3 (between the 2020-01-01 and 2020-07-01 rows),
# shifting all subsequent rows down by one, and populate it with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3, pushing rows 3..22 down to 4..23.
$ws.Range("A3:H3").Insert("xlShiftDown")

# The inserted row picked up a blended/new style; re-apply the same
# formatting used by the other label cells in column A.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# Populate the newly inserted row with the matched-error data for the
# 2020-04-01 quarter.
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value = 2.174397541324862
$ws.Range("C3").Value = 9.643547872324863
$ws.Range("D3").Value = -8.527713777675137
$ws.Range("E3").Value = -0.3880427776751375
$ws.Range("F3").Value = 1.759528090324862
$ws.Range("G3").Value = -1.680501777675137
$ws.Range("H3").Value = -1.741675777675137
